$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 467.5
$ws.Range("I20").Value = 467.5
$ws.Range("K20").Value = 467.5
$ws.Range("M20").Value = -237.5
$ws.Range("H35").Value = 467.5
$ws.Range("I35").Value = 467.5
$ws.Range("K35").Value = 467.5
$ws.Range("M35").Value = -88.5
$ws.Range("H88").Value = 1342.4615
$ws.Range("I88").Value = 267.6
$ws.Range("K88").Value = 267.6
$ws.Range("M88").Value = 138.4
$ws.Range("H91").Value = 1342.4615
$ws.Range("I91").Value = 267.6
$ws.Range("K91").Value = 267.6
$ws.Range("M91").Value = 1136.4
$ws.Range("H111").Value = 4005.2
$ws.Range("J111").Value = 2999
$ws.Range("L111").Value = 8997
$ws.Range("N111").Value = -15131
$ws.Range("H129").Value = 200000
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").Value = ""
$ws.Range("H132").Value = 1383.5217
$ws.Range("I132").Value = 1383.5217
$ws.Range("K132").Value = 4150.5651
$ws.Range("M132").Value = -1620.5651
$ws.Range("H137").Value = 2696.4
$ws.Range("I137").Value = 2045
$ws.Range("K137").Value = 6135
$ws.Range("M137").Value = -3585

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7490.143
$ws.Range("I32").Value = 5647
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 5647
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -5360
$ws.Range("N32").Value = -25574
$ws.Range("H45").Value = 1831.7142
$ws.Range("I45").Value = 1818.7693
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1818.7693
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1441.7693
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 1391.9166
$ws.Range("I61").Value = 1391.9166
$ws.Range("K61").Value = 1391.9166
$ws.Range("M61").Value = -1179.9166
$ws.Range("H74").Value = 1790.5
$ws.Range("I74").Value = 1001.8889
$ws.Range("J74").Value = 8888
$ws.Range("K74").Value = 1001.8889
$ws.Range("L74").Value = 8888
$ws.Range("M74").Value = -127.8889
$ws.Range("N74").Value = -10636
$ws.Range("H77").Value = 1790.5
$ws.Range("I77").Value = 1001.8889
$ws.Range("J77").Value = 8888
$ws.Range("K77").Value = 5009.444500000001
$ws.Range("L77").Value = 44440
$ws.Range("M77").Value = -641.4445000000005
$ws.Range("N77").Value = -53176
$ws.Range("H88").Value = 400
$ws.Range("I88").Value = 400
$ws.Range("K88").Value = 400
$ws.Range("M88").Value = 6
$ws.Range("H91").Value = 400
$ws.Range("I91").Value = 400
$ws.Range("K91").Value = 400
$ws.Range("M91").Value = 1004
$ws.Range("H97").Value = 1026.4166
$ws.Range("J97").Value = 892.3333
$ws.Range("L97").Value = 892.3333
$ws.Range("N97").Value = -1884.3333
$ws.Range("H132").Value = 1426.4445
$ws.Range("I132").Value = 1356.9286
$ws.Range("K132").Value = 4070.7858
$ws.Range("M132").Value = -1540.7858
$ws.Range("H136").Value = 1391.9166
$ws.Range("I136").Value = 1391.9166
$ws.Range("K136").Value = 4175.7498
$ws.Range("M136").Value = -1625.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1496
$ws.Range("I86").Value = 1493.4
$ws.Range("J86").Value = 1499.25
$ws.Range("K86").Value = 1493.4
$ws.Range("L86").Value = 1499.25
$ws.Range("M86").Value = -370.4000000000001
$ws.Range("N86").Value = -3745.25
$ws.Range("H89").Value = 1496
$ws.Range("I89").Value = 1493.4
$ws.Range("J89").Value = 1499.25
$ws.Range("K89").Value = 7467
$ws.Range("L89").Value = 7496.25
$ws.Range("M89").Value = -1851
$ws.Range("N89").Value = -18728.25
$ws.Range("H107").Value = 710.55
$ws.Range("I107").Value = 706.2222
$ws.Range("J107").Value = 749.5
$ws.Range("K107").Value = 706.2222
$ws.Range("L107").Value = 749.5
$ws.Range("M107").Value = 1213.7778
$ws.Range("N107").Value = -4589.5
$ws.Range("H134").Value = 2218.5715
$ws.Range("I134").Value = 2854
$ws.Range("K134").Value = 8562
$ws.Range("M134").Value = -6027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2776.889
$ws.Range("I58").Value = 1101.25
$ws.Range("J58").Value = 6128.1665
$ws.Range("K58").Value = 1101.25
$ws.Range("L58").Value = 6128.1665
$ws.Range("M58").Value = -898.25
$ws.Range("N58").Value = -6534.1665
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H105").Value = 1747.909
$ws.Range("I105").Value = 818.2857
$ws.Range("J105").Value = 3374.75
$ws.Range("K105").Value = 818.2857
$ws.Range("L105").Value = 3374.75
$ws.Range("M105").Value = 928.7143
$ws.Range("N105").Value = -6868.75
$ws.Range("H106").Value = 21402.4
$ws.Range("J106").Value = 21402.4
$ws.Range("L106").Value = 21402.4
$ws.Range("N106").Value = -23926.4
$ws.Range("H107").Value = 660.86365
$ws.Range("I107").Value = 485.2
$ws.Range("K107").Value = 485.2
$ws.Range("M107").Value = 1434.8
$ws.Range("H132").Value = 2137.6956
$ws.Range("I132").Value = 1291.125
$ws.Range("J132").Value = 4072.7144
$ws.Range("K132").Value = 3873.375
$ws.Range("L132").Value = 12218.1432
$ws.Range("M132").Value = -1343.375
$ws.Range("N132").Value = -17278.1432
$ws.Range("H134").Value = 3058.2942
$ws.Range("I134").Value = 2561.4666
$ws.Range("J134").Value = 6784.5
$ws.Range("K134").Value = 7684.399800000001
$ws.Range("L134").Value = 20353.5
$ws.Range("M134").Value = -5149.399800000001
$ws.Range("N134").Value = -25423.5
$ws.Range("H136").Value = 2776.889
$ws.Range("I136").Value = 1101.25
$ws.Range("J136").Value = 6128.1665
$ws.Range("K136").Value = 3303.75
$ws.Range("L136").Value = 18384.4995
$ws.Range("M136").Value = -753.75
$ws.Range("N136").Value = -23484.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 14285773
$ws.Range("I7").Value = 33333378
$ws.Range("K7").Value = 100000134
$ws.Range("M7").Value = -100000022
$ws.Range("H48").Value = 550
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 550
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 1650
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = -2150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1571.3572
$ws.Range("I107").Value = 2999.75
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2999.75
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -1079.75
$ws.Range("N107").Value = -4840
$ws.Range("H113").Value = 4224.222
$ws.Range("I113").Value = 3254.5
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3254.5
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1084.5
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3605.3845
$ws.Range("I61").Value = 4639
$ws.Range("J61").Value = 2399.5
$ws.Range("K61").Value = 4639
$ws.Range("L61").Value = 2399.5
$ws.Range("M61").Value = -4437
$ws.Range("N61").Value = -2803.5
$ws.Range("H100").Value = 3959.9
$ws.Range("I100").Value = 1779.8
$ws.Range("J100").Value = 6140
$ws.Range("K100").Value = 1779.8
$ws.Range("L100").Value = 6140
$ws.Range("M100").Value = -1238.8
$ws.Range("N100").Value = -7222
$ws.Range("H110").Value = 59999
$ws.Range("J110").Value = 59999
$ws.Range("L110").Value = 59999
$ws.Range("N110").Value = -68179
$ws.Range("H113").Value = 3605.3845
$ws.Range("I113").Value = 4639
$ws.Range("J113").Value = 2399.5
$ws.Range("K113").Value = 4639
$ws.Range("L113").Value = 2399.5
$ws.Range("M113").Value = -2469
$ws.Range("N113").Value = -6739.5
$ws.Range("H132").Value = 3624.6924
$ws.Range("I132").Value = 2802.6875
$ws.Range("J132").Value = 4939.9
$ws.Range("K132").Value = 8408.0625
$ws.Range("L132").Value = 14819.7
$ws.Range("M132").Value = -5878.0625
$ws.Range("N132").Value = -19879.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 38019
$ws.Range("I51").Value = 34999.5
$ws.Range("J51").Value = 41038.5
$ws.Range("K51").Value = 34999.5
$ws.Range("L51").Value = 41038.5
$ws.Range("M51").Value = -34489.5
$ws.Range("N51").Value = -42058.5
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
$ws.Range("H132").Value = 1834.4375
$ws.Range("I132").Value = 1637.25
$ws.Range("K132").Value = 4911.75
$ws.Range("M132").Value = -2381.75
$ws.Range("H136").Value = 2684.682
$ws.Range("I136").Value = 849
$ws.Range("J136").Value = 4520.364
$ws.Range("K136").Value = 2547
$ws.Range("L136").Value = 13561.092
$ws.Range("M136").Value = 3
$ws.Range("N136").Value = -18661.092
